$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 0.7
$ws.Range("B17").Value = 2985065.736644
$ws.Range("C17").Value = 21057.645339
$ws.Range("D17").Value = 2964008.091305
$ws.Range("E17").Value = 13607.360941
$ws.Range("F17").Value = 1290516.054363
$ws.Range("G17").Value = 14272.357413
$ws.Range("H17").Value = 1276243.69695
$ws.Range("I17").Value = 15537.91470366667
